# Commit: "update 20220708 18:44 by xhx"
# Sets the "variance" column (F) to 0 for a specific set of rule rows
# (previously 2), and moves the active sheet selection to I39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column F ("variance") value changes from 2 -> 0.
$rows = @(2,3,4,5,13,14,15,16,17,20,21,22,23,24,25,28,30,31,32,33,34,35,39,40,41,42,43,44,45)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = 0
}

# Update the active cell / selection on the sheet.
$ws.Range("I39").Select() | Out-Null
